$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Agosto de 2020 a las 23:28"

# Row 4
$ws.Range("B4").Value = 6169205
$ws.Range("C4").Value = 29835
$ws.Range("D4").Value = 3421866
$ws.Range("E4").Value = 2560145
$ws.Range("G4").Value = 339
$ws.Range("H4").Value = 187194

# Row 8
$ws.Range("B8").Value = 647166
$ws.Range("C8").Value = 7731
$ws.Range("D8").Value = 455457
$ws.Range("E8").Value = 162921
$ws.Range("G8").Value = 181
$ws.Range("H8").Value = 28788

# Row 9
$ws.Range("B9").Value = 625056
$ws.Range("C9").Value = 2505
$ws.Range("D9").Value = 538604
$ws.Range("E9").Value = 72424
$ws.Range("G9").Value = 47
$ws.Range("H9").Value = 14028

# Row 23
$ws.Range("B23").Value = 243295
$ws.Range("C23").Value = 470
$ws.Range("E23").Value = 16447

# Row 54
$ws.Range("B54").Value = 51574
$ws.Range("C54").Value = 183
$ws.Range("D54").Value = 48654
$ws.Range("E54").Value = 2731

# Row 189
$ws.Range("B189").Value = 173
$ws.Range("C189").Value = 3
$ws.Range("D189").Value = 145
$ws.Range("E189").Value = 21

# Row 218
$ws.Range("D218").Value = 4
$ws.Range("E218").Value = 1
